$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 79500
$ws.Range("B3").Value = 91829
$ws.Range("B4").Value = 80350
$ws.Range("B5").Value = 80349
